$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the data table with just the header row labels "sku" / "urun_adı",
# clearing out all the former data rows (2-9).
$ws.Range("A2:B9").ClearContents() | Out-Null
$ws.Range("A1").Value = "sku"
$ws.Range("B1").Value = "urun_adı"

# Move the saved cell selection to D11 (previously B11).
$ws.Range("D11").Select() | Out-Null
